# Append the new data row (row 81) reported by the automatic PEBCOM map update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 81

# Columns that must be stored as literal text (matches how every other row
# in this sheet stores "numeric-looking" values like Caso/OT/Comuna/dates).
$textCols = @(1, 2, 3, 4, 5, 6, 7, 8, 10, 11, 12, 15, 16)
foreach ($c in $textCols) {
    $ws.Cells.Item($row, $c).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value  = "-541"
$ws.Cells.Item($row, 2).Value  = "8/1/2025"
$ws.Cells.Item($row, 3).Value  = "AYACUCHO 241"
$ws.Cells.Item($row, 4).Value  = "3"
$ws.Cells.Item($row, 5).Value  = "808663880"
$ws.Cells.Item($row, 6).Value  = "PEBCOM"
$ws.Cells.Item($row, 7).Value  = "Pendiente"
$ws.Cells.Item($row, 8).Value  = "Colocar columna para pedir traspaso de nodo"
$ws.Cells.Item($row, 9).Value  = 1
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Nodo Teco"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.395015
$ws.Cells.Item($row, 14).Value = -34.606755
$ws.Cells.Item($row, 15).Value = "Almagro"
$ws.Cells.Item($row, 16).Value = "Capital Sur"

# Drop the temporary "@" text format again so the new row keeps the same
# (default/no explicit style) formatting as the rest of the sheet.
foreach ($c in $textCols) {
    $ws.Cells.Item($row, $c).ClearFormats()
}
